# sizework20240216.xlsx update
# - removing old WS Q param conversion; reuse the div param better
# - slope / filter captions
# - biquad filter is now an option in Maj7
# - moog filter is now disableable (safer for size); biquad + sat kinda replaces it.
#
# The shared-strings "w/o maj7", "Sum", "Average", "Running Total", "Count" were
# already unused by any cell, so they simply drop out of the saved sharedStrings
# table once the workbook is re-saved (Excel only persists strings that are still
# referenced) - no explicit action required for that part of the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Row 61: fill in the previously-empty "diff note" column F (same text it always
# had logically, the shared-string table just renumbers around it automatically).
$ws.Cells.Item(61, 6).Value = "at least it was a tiny change"

# Row 62 (new data point): biquad option / onepole readded / moog removed.
$ws.Cells.Item(62, 1).Value = "adding biquad option, readding onepole, removing moog"
$ws.Cells.Item(62, 2).Value = 20172

# Row 63 (new data point): removing obsolete param types.
$ws.Cells.Item(63, 1).Value = "also removing obsolete param types"
$ws.Cells.Item(63, 2).Value = 20172
$ws.Cells.Item(63, 6).Value = "at this point I have saved well over 1kb, though added"

# Row 64 (new data point): removing comp & sat.
$ws.Cells.Item(64, 1).Value = "removing comp & sat"
$ws.Cells.Item(64, 2).Value = 18812
$ws.Cells.Item(64, 6).Value = "1.4kb via sat & comp. but you can argue I'll save this via directives. So I'm good."

# Columns C/D/E for rows 62-64 (and the downstream total row 65) are driven by the
# existing shared formulas already present in the sheet (C = B(row)-B(row-1),
# E = B(row)-D(row)), so they recompute on their own once B62:B64 are populated.

# Restore the view's active selection to match the latest edit location.
$ws.Range("F65").Select()
